$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.522.39"
$ws.Range("E2").Value = "  +1.30%  "

$ws.Range("D3").Value = "2.983.51"
$ws.Range("E3").Value = "  +2.65%  "

$ws.Range("E4").Value = "  +0.44%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "381.21"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.55%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "104.59"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.10%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.545"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.05%  "

$ws.Range("E8").Value = "  +0.21%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.597"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.63%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "37.31"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.75%  "

$ws.Range("E11").Value = "  +0.23%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0848"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.98%  "

$ws.Range("D13").Value = "3.454.54"
$ws.Range("E13").Value = "  +3.07%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.42"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.02%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.58"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.09%  "

$ws.Range("D16").Value = "2.988.57"
$ws.Range("E16").Value = "  +2.93%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.977"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +5.92%  "

$ws.Range("D18").Value = "51.530.83"
$ws.Range("E18").Value = "  +1.94%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.32"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.78%  "

$ws.Range("E20").Value = "  +4.25%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.98"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.90%  "

$ws.Range("D22").Value = "0.0₃0965"
$ws.Range("E22").Value = "  +2.68%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.42"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.19%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "262.97"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.04%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.94"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +10.26%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.28"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +17.34%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.70"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +23.05%  "

$ws.Range("E28").Value = "  +14.08%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.171"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.75%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "26.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.90%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.999"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.01%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "9.89"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.50%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "35.01"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.44%  "

$ws.Range("E34").Value = "  -1.58%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "51.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.23%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0453"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +7.97%  "

$ws.Range("E37").Value = "  +0.55%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.05"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.56%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "17.22"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.40%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.59"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.06%  "

$ws.Range("E41").Value = "  +0.76%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.116"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.75%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "125.52"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.59%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "21.79"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.78%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.282"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +19.90%  "

$ws.Range("E46").Value = "  -1.47%  "

$ws.Range("E47").Value = "  +2.46%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.28"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.88%  "

$ws.Range("D49").Value = "2.032.93"
$ws.Range("E49").Value = "  +0.90%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0331"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +7.70%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "58.44"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.15%  "
